# Update countries & provincias Spain
# Daily data refresh: reorder a few tied/updated countries and refresh
# case numbers, plus bump the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

function Set-Row($row, $values) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

# --- Timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 24 de Mayo de 2020 a las 21:35"

# --- Straight numeric refresh (country order unchanged) ---
# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1677614
$ws.Cells.Item(4, 3).Value = 10786
$ws.Cells.Item(4, 5).Value = 1127391
$ws.Cells.Item(4, 7).Value = 364
$ws.Cells.Item(4, 8).Value = 99047

# Row 34: Sudafrica
$ws.Cells.Item(34, 4).Value = 11100
$ws.Cells.Item(34, 5).Value = 11054

# Row 58: Oman
$ws.Cells.Item(58, 5).Value = 5800
$ws.Cells.Item(58, 7).Value = 1
$ws.Cells.Item(58, 8).Value = 37

# --- Rows whose country swapped place with its neighbour, with new data ---
# Row 140 now shows Togo (was Cabo Verde) with freshly updated figures.
Set-Row 140 @("Togo", 381, 8, 141, 228, 0, 0, 12)
# Row 141 now shows Cabo Verde (was Togo), carrying the old Cabo Verde figures.
Set-Row 141 @("Cabo Verde", 380, 9, 155, 222, 0, 0, 3)

# Row 198 now shows Nueva Caledonia (was Santa Lucia); figures tied, unchanged.
Set-Row 198 @("Nueva Caledonia", 18, 0, 18, 0, 0, 0, 0)
# Row 199 now shows Santa Lucia (was Nueva Caledonia); figures tied, unchanged.
Set-Row 199 @("Santa Lucia", 18, 0, 18, 0, 0, 0, 0)

# Row 209 now shows Seychelles (was Groenlandia); figures tied, unchanged.
Set-Row 209 @("Seychelles", 11, 0, 11, 0, 0, 0, 0)
# Row 210 now shows Groenlandia (was Seychelles); figures tied, unchanged.
Set-Row 210 @("Groenlandia", 11, 0, 11, 0, 0, 0, 0)

# Row 214 now shows Bonaire, San Eustaquio y Saba (was Sahara Occidental).
Set-Row 214 @("Bonaire, San Eustaquio y Saba", 6, 0, 6, 0, 0, 0, 0)
# Row 215 now shows San Bartolome (was Bonaire, San Eustaquio y Saba).
Set-Row 215 @("San Bartolome", 6, 0, 6, 0, 0, 0, 0)
# Row 216 now shows Sahara Occidental (was San Bartolome).
Set-Row 216 @("Sahara Occidental", 6, 0, 6, 0, 0, 0, 0)
